$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Price (D) and Volume(1h) (E) figures pulled from the crypto feed.
# Values in column D that would otherwise be auto-parsed as numbers by Excel
# are entered with a leading apostrophe so they stay literal text (matching the
# source feed formatting, e.g. trailing zeros such as "1.00" or "7.30").

# Row 2
$ws.Range("D2").Value = "58.016.18"
$ws.Range("E2").Value = "  +0.50%  "
# Row 3
$ws.Range("D3").Value = "3.135.12"
$ws.Range("E3").Value = "  +0.27%  "
# Row 4
$ws.Range("E4").Value = "  +0.04%  "
# Row 5
$ws.Range("D5").Value = "'532.54"
$ws.Range("E5").Value = "  +0.99%  "
# Row 6
$ws.Range("D6").Value = "'138.12"
$ws.Range("E6").Value = "  -0.26%  "
# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.03%  "
# Row 8
$ws.Range("D8").Value = "3.133.54"
$ws.Range("E8").Value = "  +0.25%  "
# Row 9
$ws.Range("D9").Value = "'0.467"
$ws.Range("E9").Value = "  +4.76%  "
# Row 10
$ws.Range("D10").Value = "'7.30"
$ws.Range("E10").Value = "  +1.89%  "
# Row 11
$ws.Range("E11").Value = "  -0.54%  "
# Row 12
$ws.Range("D12").Value = "'0.412"
$ws.Range("E12").Value = "  +4.02%  "
# Row 13
$ws.Range("D13").Value = "3.670.20"
$ws.Range("E13").Value = "  +0.20%  "
# Row 14
$ws.Range("E14").Value = "  +1.41%  "
# Row 15
$ws.Range("D15").Value = "'25.69"
$ws.Range("E15").Value = "  +0.55%  "
# Row 16
$ws.Range("E16").Value = "  -0.38%  "
# Row 17
$ws.Range("D17").Value = "58.091.49"
$ws.Range("E17").Value = "  +0.51%  "
# Row 18
$ws.Range("D18").Value = "3.135.35"
$ws.Range("E18").Value = "  +0.30%  "
# Row 19
$ws.Range("D19").Value = "'6.02"
$ws.Range("E19").Value = "  +0.15%  "
# Row 20
$ws.Range("D20").Value = "'12.69"
$ws.Range("E20").Value = "  -0.66%  "
# Row 21
$ws.Range("E21").Value = "  +2.47%  "
# Row 22
$ws.Range("D22").Value = "'357.45"
$ws.Range("E22").Value = "  +1.37%  "
# Row 23
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.42%  "
# Row 24
$ws.Range("D24").Value = "'69.13"
$ws.Range("E24").Value = "  +1.06%  "
# Row 25
$ws.Range("E25").Value = "  -0.34%  "
# Row 26
$ws.Range("E26").Value = "  -1.41%  "
# Row 27
$ws.Range("E27").Value = "  +0.10%  "
# Row 28
$ws.Range("D28").Value = "0.0₃0875"
$ws.Range("E28").Value = "  -4.40%  "
# Row 29
$ws.Range("D29").Value = "'7.30"
$ws.Range("E29").Value = "  -2.29%  "
# Row 30
$ws.Range("D30").Value = "'6.17"
$ws.Range("E30").Value = "  -0.23%  "
# Row 31
$ws.Range("E31").Value = "  -0.32%  "
# Row 32
$ws.Range("D32").Value = "'21.46"
$ws.Range("E32").Value = "  +1.63%  "
# Row 33
$ws.Range("D33").Value = "'5.01"
$ws.Range("E33").Value = "  +1.27%  "
# Row 34
$ws.Range("D34").Value = "'1.14"
$ws.Range("E34").Value = "  -3.53%  "
# Row 35
$ws.Range("D35").Value = "'158.81"
$ws.Range("E35").Value = "  +0.72%  "
# Row 36
$ws.Range("D36").Value = "'6.07"
$ws.Range("E36").Value = "  -1.66%  "
# Row 37
$ws.Range("D37").Value = "'25.83"
$ws.Range("E37").Value = "  -1.54%  "
# Row 38
$ws.Range("D38").Value = "'1.26"
$ws.Range("E38").Value = "  -0.99%  "
# Row 39
$ws.Range("D39").Value = "'1.69"
$ws.Range("E39").Value = "  +3.85%  "
# Row 40
$ws.Range("E40").Value = "  +0.18%  "
# Row 41
$ws.Range("D41").Value = "2.505.20"
$ws.Range("E41").Value = "  +8.20%  "
# Row 42
$ws.Range("E42").Value = "  -0.17%  "
# Row 43
$ws.Range("E43").Value = "  -4.22%  "
# Row 44
$ws.Range("D44").Value = "'37.54"
$ws.Range("E44").Value = "  +2.81%  "
# Row 45
$ws.Range("D45").Value = "3.175.57"
$ws.Range("E45").Value = "  +0.36%  "
# Row 48
$ws.Range("D48").Value = "'0.979"
$ws.Range("E48").Value = "  +0.53%  "
# Row 49
$ws.Range("E49").Value = "  -0.10%  "
# Row 50
$ws.Range("D50").Value = "'19.79"
$ws.Range("E50").Value = "  -3.40%  "
# Row 51
$ws.Range("D51").Value = "'0.735"
$ws.Range("E51").Value = "  -4.10%  "

# Rows 46/47 swap rank order: VeChain moves above FirstDigitalUSD.
# (Column A rank-index values are unchanged; only Coin/Link/Price/Volume move.)
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0269"
$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.03%  "
